$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    67477.66905413473,
    64122.11951149303,
    62574.13871349033,
    60973.49679513348,
    60879.16604411871,
    63882.90556188919,
    63938.28624512946,
    70529.41023952424,
    84189.47718710327,
    90508.1221890134,
    92875.42047626144,
    94727.00078499186,
    91501.91594179011,
    94896.521380304,
    95445.06239869088,
    94088.23803787978,
    89152.05982478207,
    83412.98811676129,
    83554.63872678358,
    80323.43696802825,
    79037.85700430113,
    76992.95007782558,
    74412.76984318219,
    70855.3773670299
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
